# AutomationFramework - final
# Adds a new "patientErrorMessageTest" worksheet with sample test data,
# fixes up the selection left on "invalidCredentialTest", and makes the
# new sheet the active/selected tab.

$wb = $excel.ActiveWorkbook

# --- 1. Tidy up the selection remembered on sheet 1 (invalidCredentialTest) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1:C1").Select() | Out-Null

# --- 2. Add the new sheet at the end of the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "patientErrorMessageTest"

# --- 3. Header row ---
$headers = @("Username", "Password", "Language", "Message", "Type", "User Id", "ExpectedValue")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws3.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 4. Data rows ---
$row2 = @("admin", "pass", "English (Indian)", "Please book an appointment", "Insurance", 2, "Please choose a patient")
for ($i = 0; $i -lt $row2.Length; $i++) {
    $ws3.Cells.Item(2, $i + 1).Value = $row2[$i]
}

$row3 = @("accountant", "accountant", "English (Indian)", "Please book an appointment", "Insurance", 2, "Please choose a patient")
for ($i = 0; $i -lt $row3.Length; $i++) {
    $ws3.Cells.Item(3, $i + 1).Value = $row3[$i]
}

# --- 5. Column widths (best-fit-ish, matching the authored sheet) ---
$ws3.Columns.Item(1).ColumnWidth = 9.166666666666666
$ws3.Columns.Item(2).ColumnWidth = 8.592447916666666
$ws3.Columns.Item(3).ColumnWidth = 14.022135416666666
$ws3.Columns.Item(4).ColumnWidth = 26.022135416666668
$ws3.Columns.Item(7).ColumnWidth = 17.592447916666668

# --- 6. Selection + active tab on the new sheet ---
$ws3.Range("C4").Select() | Out-Null
$ws3.Activate() | Out-Null
